# regen sval data to filter save games
# Updates the per-row stat values (TB, d2S, K, IP, sum) on Sheet1 rows 2-8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New values: row -> @{ B=..; C=..; D=..; E=..; G=.. }
$newValues = @{
    2 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.8054896365839992; E = 0.496779210170732;  G = 6.201049113329182 }
    3 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.1575252929769615; E = 8.660232485948974;  G = 13.71653804550039 }
    4 = @{ B = 0.6753301551942219; C = 0.3127903958511391; D = 0.1575252929769615; E = 0.496779210170732;  G = 1.642425054193055 }
    5 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.1575252929769615; E = 0.496779210170732;  G = 5.553084769722144 }
    6 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.1575252929769615; E = 0.496779210170732;  G = 5.553084769722144 }
    7 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.1575252929769615; E = 0.496779210170732;  G = 5.553084769722144 }
    8 = @{ B = 3.230985683306322;  C = 0.3127903958511391; D = 0.8054896365839992; E = 0.496779210170732;  G = 4.846044925912192 }
}

foreach ($row in $newValues.Keys) {
    $cols = $newValues[$row]
    $ws.Range("B$row").Value = $cols.B
    $ws.Range("C$row").Value = $cols.C
    $ws.Range("D$row").Value = $cols.D
    $ws.Range("E$row").Value = $cols.E
    $ws.Range("G$row").Value = $cols.G
}
